$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: add E6:H6 (effort-burndown style values + a "/" marker) ---
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 3
$ws.Cells.Item(6, 7).Value = 2
$ws.Cells.Item(6, 8).Value = "/"

# Match formatting: E6 like D6 (wrap-style), F6/G6/H6 like A7 (task-header style)
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(6, 5).PasteSpecial(-4122)

$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(6, 6).PasteSpecial(-4122)
$ws.Cells.Item(6, 7).PasteSpecial(-4122)
$ws.Cells.Item(6, 8).PasteSpecial(-4122)

# --- Row 7: add E7:F7 ---
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = "/"

$ws.Cells.Item(7, 4).Copy()
$ws.Cells.Item(7, 5).PasteSpecial(-4122)

$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(7, 6).PasteSpecial(-4122)

# --- Row 8: add E8:G8 ---
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = "/"

$ws.Cells.Item(8, 4).Copy()
$ws.Cells.Item(8, 5).PasteSpecial(-4122)
$ws.Cells.Item(8, 6).PasteSpecial(-4122)

$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(8, 7).PasteSpecial(-4122)
